# Apply changes described in commit "Add data for 2022-03-15"
#  - Rename worksheet title / header from "2022-03-06" to "2022-03-07"
#  - Update March (row 4) "2022 through" column value 33 -> 36
#  - Update Total (row 14) "2022 through" column value 334 -> 337

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the sheet (tab) name
$ws.Name = "Through 2022-03-07"

# Update the header cell text for the "2022 (through 03-06)" column -> 03-07
$ws.Range("I1").Value = "2022 (through 03-07)"

# Update the March value in the 2022 column
$ws.Range("I4").Value = 36

# Update the Total value in the 2022 column
$ws.Range("I14").Value = 337
